# Parque_Vehicular_AF.xlsx update
# The sheet holding the monthly fleet-size series ("C_11") is renamed to
# "C_18" to reflect the newer period label used for this upload. Excel
# automatically keeps the sheet-scoped _xlnm._FilterDatabase defined name
# (and any other sheet-qualified references) pointed at the renamed sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Name = "C_18"
